$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) — J1, K1
$ws.Range("J1").Value = "Column contains texts in incorrect number format"
$ws.Range("K1").Value = "Column contains both text and numeric"

# Row 2 new cells
$ws.Range("J2").Value = "401491.00.00"
$ws.Range("K2").Value = 1234121
$ws.Range("K2").NumberFormat = "General"

# Row 3 new cells
$ws.Range("J3").Value = "3100310.00.00"
$ws.Range("K3").Value = "Lorem Ipsum"

# Column widths
$ws.Columns.Item(9).ColumnWidth = 68.7142857142857
$ws.Columns.Item(10).ColumnWidth = 52.2857142857143
$ws.Columns.Item(11).ColumnWidth = 49.5714285714286

# Selection moves to K1
$ws.Range("K1").Select()
